$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final dataset for rows 2..20 (column A = numeric id, column B = UPN string)
$data = @(
    @(7,  "A931383810034"),
    @(24, "B931100609017"),
    @(27, "Z931412020040"),
    @(46, "V931325309014"),
    @(61, "C931412020042"),
    @(11, "W931101109061"),
    @(18, "C673336808019"),
    @(25, "L931100509013"),
    @(34, "X931412020027"),
    @(57, "D931100609028"),
    @(47, "X931325208068"),
    @(44, "U931383908031"),
    @(45, "Q931321008053"),
    @(49, "Q931101109046"),
    @(2,  "R928218115049"),
    @(40, "L931412020028"),
    @(10, "Q879418719002"),
    @(41, "N931383610018"),
    @(38, "M931252916068")
)

# Reference cell that already carries the formatting (bold, border, centered)
# used throughout column A.
$ws.Range("A2").Copy()

$rowIndex = 2
foreach ($entry in $data) {
    # For rows beyond the sheet's original extent, first clone the A-column
    # formatting so new rows look consistent with the existing ones.
    if ($rowIndex -gt 12) {
        $ws.Range("A$rowIndex").PasteSpecial(-4122)  # xlPasteFormats
    }

    $ws.Cells.Item($rowIndex, 1).Value = $entry[0]
    $ws.Cells.Item($rowIndex, 2).Value = $entry[1]

    $rowIndex++
}

$excel.CutCopyMode = 0
